$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Add "NA" values in column E (duplicate_image_filename) for rows 2 through 21
$ws.Range("E2:E21").Value = "NA"
